$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: update abstract text (row 2)
$ws.Range("D2").Value = 'Abstract Introduction The mechanisms for the observed disparities in diabetes‐related amputation are poorly understood and could be related to access for diabetic foot ulceration (DFU) care. This qualitative study aimed to understand patients'' personal experiences navigating the healthcare system and the barriers they faced. Methods Fifteen semi‐structured interviews were conducted over the phone between June 2020 to February 2021. Participants with DFUs were recruited from a tertiary referral center in Southern Arizona. The interviews were audio‐recorded and analyzed according to the NIMHD Research Framework, focusing on the health care system domain. Results Among the 15 participants included in the study, the mean age was 52.4 years (66.7% male), 66.7% was from minority racial groups, and 73.3% was Medicaid or Indian Health Service beneficiaries. Participants frequently reported barriers at various levels of the healthcare system. On the individual level, themes that arose included health literacy and inadequate insurance coverage resulting in financial strain. On the interpersonal level, participants complained of fragmented relationships with providers and experienced challenges in making follow‐up appointments. On the community level, participants reported struggles with medical equipment. On the societal level, participants also noted insufficient preventative foot care and education before DFU onset, and many respondents experienced initial misdiagnoses and delays in receiving care. Conclusions Patients with DFUs face significant barriers in accessing medical care at many levels in the healthcare system and beyond. These data highlight opportunities to address the effects of diabetic foot complications and the inequitable burden of inadequately managed diabetic foot care.'

# E2: update publication_date (row 2) - force text via leading apostrophe so Excel
# does not reinterpret the date-like string as a serial date number
$ws.Range("E2").Value = '''2022-01-01'

# AF2: update related_works list (row 2)
$ws.Range("AF2").Value = 'c("https://openalex.org/W2107218281", "https://openalex.org/W2325887025", "https://openalex.org/W2138721938", "https://openalex.org/W2992118406", "https://openalex.org/W3113711756", "https://openalex.org/W4388073528", "https://openalex.org/W2375549631", "https://openalex.org/W2347618812", "https://openalex.org/W3031691025", "https://openalex.org/W2468650808")'

# D4: fix curly apostrophe -> straight apostrophe in "residents " (row 4)
$ws.Range("D4").Value = 'Perceptions of community can play an important role in determining health and well-being. We know little, however, about residents'' perceptions of community safety in the Southwestern borderlands, an area frequently portrayed as plagued by disorder. The qualitative aim of this community-based participatory research study was to explore the perceptions of Mexican-origin border residents about their communities in southern Yuma County, Arizona. Our team of University of Arizona researchers and staff from Campesinos Sin Fronteras, a grassroots farmworker support agency in Yuma County, Arizona, developed a bilingual interview guide and recruited participants through radio adds, flyers, and cold calls among existing agency clientele. Thirty individual interviews with participants of Mexican origin who live in and/or work in rural Yuma County were conducted remotely in 2021. Participants overwhelmingly perceived their communities as both calm and safe. While some participants mentioned safety concerns, the vast majority described high levels of personal security and credited both neighbors and police for ensuring local safety. These perceptions were stated in direct contrast to those across the border, where participants had positive familial and cultural ties but negative perceptions regarding widespread violence. In conclusion, we argue that to understand environmental factors affecting health and well-being in Mexican immigrant populations, it is critical to examine the role of binational external referents that color community perceptions.'

# AF4: update related_works list (row 4)
$ws.Range("AF4").Value = 'c("https://openalex.org/W2748952813", "https://openalex.org/W2891252757", "https://openalex.org/W2005441443", "https://openalex.org/W2273917609", "https://openalex.org/W136810046", "https://openalex.org/W3123059865", "https://openalex.org/W2587158935", "https://openalex.org/W3126052823", "https://openalex.org/W3124033816", "https://openalex.org/W3212206169")'

# AF5: update related_works list (row 5) - reordered / replaced one id
$ws.Range("AF5").Value = 'c("https://openalex.org/W4251033897", "https://openalex.org/W2084603631", "https://openalex.org/W2028171831", "https://openalex.org/W4379418727", "https://openalex.org/W4379419038", "https://openalex.org/W2282792614", "https://openalex.org/W4281771988", "https://openalex.org/W2115024886", "https://openalex.org/W2050646317", "https://openalex.org/W2112228454")'

